# reportRetail.xlsx — add "Month" and "Sales Manager" breakdowns
#
# Adds two new columns (D: Month, E: Sales Manager) to the existing
# retail-category table and appends two small concatenated result
# blocks (a per-month sales summary and a per-sales-manager sales
# summary) underneath the original 36 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (D1/E1), formatted like the existing header row ---
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)   # xlPasteFormats — keep header look (bold/border/center)

$ws.Cells.Item(1, 4).Value = "Month"
$ws.Cells.Item(1, 5).Value = "Sales Manager"

# --- Appended "Month" summary block (rows 38-42) ---
$monthRows = @(
    @(2954149, 17.38032727376488, "April"),
    @(3597561, 21.16575283350057, "Feburary"),
    @(3610492, 21.24183058447963, "January"),
    @(3422543, 20.13605862416997, "March"),
    @(3412340, 20.07603068408495, "May")
)

$r = 38
foreach ($row in $monthRows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

# --- Appended "Sales Manager" summary block (rows 43-50) ---
$managerRows = @(
    @(29181,   0.1716823796550997,  "Chen Cho"),
    @(6750743, 39.71706324937482,   "Dominique Kai"),
    @(2544106, 14.96789596569059,   "Donald Ducker"),
    @(2540790, 14.94838673807891,   "Jane Maria"),
    @(1297791, 7.635373947944603,   "John Doe"),
    @(355958,  2.094229687031629,   "Marc Jensen"),
    @(515025,  3.030078392853834,   "Martin Miller"),
    @(2963491, 17.43528963937052,   "Vicky Dullo")
)

foreach ($row in $managerRows) {
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $r = $r + 1
}
